# Trade #72 closed at 2026-02-17 15:48:38 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets of the live trading results workbook to reflect the newly closed
# trade #72.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.25              # Current Capital
$summary.Range("B4").Value = 0.24                 # Total P&L $
$summary.Range("B5").Value = 0.07000000000000001  # Total P&L %
$summary.Range("B6").Value = 72                   # Total Trades
$summary.Range("B7").Value = 23                   # Winning Trades
$summary.Range("B9").Value = 31.94                # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.25                # Capital
$status.Range("D4").Value = 72                    # Trades
$status.Range("E4").Value = 0.24                  # P&L $
$status.Range("F4").Value = 0.25                  # P&L %
$status.Range("G4").Value = 31.94                 # Win Rate %

# ---------------------------------------------------------------------
# Helper that appends the closed-trade row #72 (sheet row 73) to a
# trades-log sheet ("All Trades" and "MarketMaking" share the schema).
# ---------------------------------------------------------------------
function Add-TradeRow72($sheet) {
    $sheet.Range("A73").Value = 72

    # Columns B and C hold plain text (a date-like string and a time-like
    # string); force text formatting so Excel does not coerce them into
    # date/time serial numbers.
    $sheet.Range("B73").NumberFormat = "@"
    $sheet.Range("B73").Value = "2026-02-17"

    $sheet.Range("C73").Value = "15:48:31"
    $sheet.Range("D73").Value = "MarketMaking"
    $sheet.Range("E73").Value = "UP"
    $sheet.Range("F73").Value = 0.72
    $sheet.Range("G73").Value = 0.8
    $sheet.Range("H73").Value = "CLOSED"
    $sheet.Range("I73").Value = 11.1111
    $sheet.Range("J73").Value = 0.08
    $sheet.Range("K73").Value = 100.25
    $sheet.Range("L73").Value = 0
    $sheet.Range("M73").Value = 0
    $sheet.Range("N73").Value = 0.6
    $sheet.Range("O73").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P73").Value = "early_exit"
    $sheet.Range("Q73").Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow72 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow72 $marketMaking

Write-Host "Trade #72 appended to All Trades / MarketMaking; Summary and Strategy Status updated."
